$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card15")

# The row that used to be the open/blank entry (row 14) becomes a filled
# historical record once a new event is logged: its empty placeholder
# cells get the literal text "nan" (matching every other completed row
# on this sheet).
$nanCols14 = @("B","C","D","E","F","G","H","I","J","K","M","P","Q","R")
foreach ($col in $nanCols14) {
    $ws.Range($col + "14").Value = "nan"
}

# New event row 15 for the service record added via Card15.
# Column A stores the card number as text (like every other row on this
# sheet), so route it through TEXT()+paste-values instead of a plain
# Value assignment (which would store it as a number).
$ws.Range("A15").Formula = "=TEXT(15,""0"")"
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4163) | Out-Null

$ws.Range("L15").Value = "1\5\2025"
$ws.Range("N15").Value = "تم سن الفلاتس +تغيير الجريده 1+سن السليندر"
$ws.Range("O15").Value = "الخبير"

$excel.CutCopyMode = $false
